$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Package info": narrow column B, update version/build metadata
# ---------------------------------------------------------------------------
$wsPkg = $wb.Worksheets.Item("Package info")
$wsPkg.Columns.Item(2).ColumnWidth = 56.833333333333336
$wsPkg.Range("B4").Value = "V1.1.2"
$wsPkg.Range("B5").Value = "DHIS2.35.3-3492688"
$wsPkg.Range("B6").Value = "20210408T081801"
$wsPkg.Range("B7").Value = "COVAC_TRACKER_V1.1.2_DHIS2.35.3-3492688_20210408T081801"

# ---------------------------------------------------------------------------
# Sheet "programs": bump the "Last updated" date
# ---------------------------------------------------------------------------
$wsProg = $wb.Worksheets.Item("programs")
$wsProg.Range("C2").Value = "2021-03-19"

# ---------------------------------------------------------------------------
# Sheet "programTrackedEntityAttributes": replace raw UIDs in column B with
# the human readable tracked entity attribute names
# ---------------------------------------------------------------------------
$wsPTEA = $wb.Worksheets.Item("programTrackedEntityAttributes")
$wsPTEA.Range("B4").Value = "First Name"
$wsPTEA.Range("B5").Value = "Surname"
$wsPTEA.Range("B6").Value = "Sex"
$wsPTEA.Range("B8").Value = "Date of birth"
$wsPTEA.Range("B10").Value = "Home Address"

# ---------------------------------------------------------------------------
# Sheet "dataElementGroups": column B values were re-ordered
# ---------------------------------------------------------------------------
$wsDEG = $wb.Worksheets.Item("dataElementGroups")
$wsDEG.Range("B3").Value = "COVAC - Underlying condition Other"
$wsDEG.Range("B4").Value = "COVAC - Dose Number"
$wsDEG.Range("B5").Value = "COVAC- Batch Number"
$wsDEG.Range("B6").Value = "COVAC - Renal Disease"
$wsDEG.Range("B7").Value = "COVAC - Malignancy"
$wsDEG.Range("B8").Value = "COVAC - Vaccine Name"
$wsDEG.Range("B9").Value = "COVAC - Pregnancy"
$wsDEG.Range("B10").Value = "COVAC - Immunodeficiency"
$wsDEG.Range("B11").Value = "COVAC Previously infected with COVID"
$wsDEG.Range("B12").Value = "COVAC Suggested date for next dose"
$wsDEG.Range("B13").Value = "COVAC - Multiple products used - Explain"
$wsDEG.Range("B14").Value = "COVAC - AEFIs present"
$wsDEG.Range("B15").Value = "COVAC - Cardiovascular Disease"
$wsDEG.Range("B16").Value = "COVAC - Vaccine Manufacturer"
$wsDEG.Range("B17").Value = "COVAC - Last Dose"
$wsDEG.Range("B18").Value = "COVAC - Pregnancy gestation"
$wsDEG.Range("B19").Value = "COVAC - Allergic reaction after first dose"
$wsDEG.Range("B20").Value = "COVAC - Chronic Lung Disease"
$wsDEG.Range("B21").Value = "COVAC - Diabetes"
$wsDEG.Range("B22").Value = "COVAC - Neurological/Neuromuscular"
$wsDEG.Range("B23").Value = "COVAC - Underlying condition"
$wsDEG.Range("B24").Value = "COVAC - Total doses"

# ---------------------------------------------------------------------------
# Sheet "trackedEntityAttributes": widen column B, grow the table from 6 to
# 11 data rows (5 new attributes) and re-establish alphabetical order by
# Name (column A). Rows keep the zebra-stripe style based on row position
# (even row -> style of row 2, odd row -> style of row 3), so rebuild every
# data row from the two template rows instead of relying on Range.Sort
# (which would carry the style along with the data instead of the row).
# ---------------------------------------------------------------------------
$wsTEA = $wb.Worksheets.Item("trackedEntityAttributes")
$wsTEA.Columns.Item(2).ColumnWidth = 21.833333333333332

# Re-stamp rows 4-12 with the correct zebra-stripe formatting (row 2's style
# for even rows, row 3's style for odd rows) before writing the final,
# alphabetically-sorted content into them.
$wsTEA.Range("A2:E2").Copy($wsTEA.Range("A4:E4"))
$wsTEA.Range("A3:E3").Copy($wsTEA.Range("A5:E5"))
$wsTEA.Range("A2:E2").Copy($wsTEA.Range("A6:E6"))
$wsTEA.Range("A3:E3").Copy($wsTEA.Range("A7:E7"))
$wsTEA.Range("A2:E2").Copy($wsTEA.Range("A8:E8"))
$wsTEA.Range("A3:E3").Copy($wsTEA.Range("A9:E9"))
$wsTEA.Range("A2:E2").Copy($wsTEA.Range("A10:E10"))
$wsTEA.Range("A3:E3").Copy($wsTEA.Range("A11:E11"))
$wsTEA.Range("A2:E2").Copy($wsTEA.Range("A12:E12"))

# Row 4: Date of birth
$wsTEA.Range("A4").Value = "Date of birth"
$wsTEA.Range("B4").Value = "patinfo_ageonsetunit"
$wsTEA.Range("C4").Value = ""
$wsTEA.Range("E4").Value = "NI0QRzJvQ0k"

# Row 5: Date of birth is estimated
$wsTEA.Range("A5").Value = "Date of birth is estimated"
$wsTEA.Range("B5").Value = ""
$wsTEA.Range("C5").Value = ""
$wsTEA.Range("E5").Value = "Z1rLc1rVHK8"

# Row 6: First Name
$wsTEA.Range("A6").Value = "First Name"
$wsTEA.Range("B6").Value = "first_name"
$wsTEA.Range("C6").Value = ""
$wsTEA.Range("E6").Value = "sB1IHYu2xQT"

# Row 7: Home Address
$wsTEA.Range("A7").Value = "Home Address"
$wsTEA.Range("B7").Value = "patinfo_resadmin0"
$wsTEA.Range("C7").Value = ""
$wsTEA.Range("E7").Value = "Xhdn49gUd52"

# Row 8: Mobile phone number
$wsTEA.Range("A8").Value = "Mobile phone number"
$wsTEA.Range("B8").Value = ""
$wsTEA.Range("C8").Value = ""
$wsTEA.Range("E8").Value = "fctSQp5nAYl"

# Row 9: National ID
$wsTEA.Range("A9").Value = "National ID"
$wsTEA.Range("B9").Value = ""
$wsTEA.Range("C9").Value = ""
$wsTEA.Range("E9").Value = "Ewi7FUfcHAD"

# Row 10: Sex
$wsTEA.Range("A10").Value = "Sex"
$wsTEA.Range("B10").Value = "patinfo_sex"
$wsTEA.Range("C10").Value = ""
$wsTEA.Range("E10").Value = "oindugucx72"

# Row 11: Surname
$wsTEA.Range("A11").Value = "Surname"
$wsTEA.Range("B11").Value = "surname"
$wsTEA.Range("C11").Value = "The patient's surname (family name)"
$wsTEA.Range("E11").Value = "ENRjVGxVL6l"

# Row 12: Unique System Identifier (EPI)
$wsTEA.Range("A12").Value = "Unique System Identifier (EPI)"
$wsTEA.Range("B12").Value = ""
$wsTEA.Range("C12").Value = "System-generated unique ID following pattern: EPI prefix + value randomly generated (#####) - Customize the length depending on the target population of your implementation"
$wsTEA.Range("E12").Value = "KSr2yTdu1AI"
